$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Amelx"
$ws.Range("C2").Value = "Lamp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.052706
$ws.Range("H2").Value = 0.158118
$ws.Range("I2").Value = 0.0131977395622021
$ws.Range("J2").Value = 0.0131977395622021
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.528487
$ws.Range("N2").Value = 85.585461
$ws.Range("O2").Value = 0.1381240089280516
$ws.Range("P2").Value = 0.1381240089280516
$ws.Range("Q2").Value = 1.503622435822
$ws.Range("R2").Value = 13.532601922398
$ws.Range("S2").Value = 0.001822924697119703
$ws.Range("T2").Value = 0.001822924697119703

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Amelx"
$ws.Range("C3").Value = "Lamp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.052706
$ws.Range("H3").Value = 0.158118
$ws.Range("I3").Value = 0.0131977395622021
$ws.Range("J3").Value = 0.0131977395622021
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 139.2310486666667
$ws.Range("N3").Value = 417.693146
$ws.Range("O3").Value = 0.6741034184216166
$ws.Range("P3").Value = 0.6741034184216166
$ws.Range("Q3").Value = 7.338311651025333
$ws.Range("R3").Value = 66.044804859228
$ws.Range("S3").Value = 0.008896641354318648
$ws.Range("T3").Value = 0.008896641354318648

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Amelx"
$ws.Range("C4").Value = "Lamp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.052706
$ws.Range("H4").Value = 0.158118
$ws.Range("I4").Value = 0.0131977395622021
$ws.Range("J4").Value = 0.0131977395622021
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 38.78302866666667
$ws.Range("N4").Value = 116.349086
$ws.Range("O4").Value = 0.1877725726503318
$ws.Range("P4").Value = 0.1877725726503318
$ws.Range("Q4").Value = 2.044098308905333
$ws.Range("R4").Value = 18.396884780148
$ws.Range("S4").Value = 0.002478173510763754
$ws.Range("T4").Value = 0.002478173510763754

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Amelx"
$ws.Range("C5").Value = "Lamp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.388571333333333
$ws.Range("H5").Value = 4.165713999999999
$ws.Range("I5").Value = 0.347702402399595
$ws.Range("J5").Value = 0.347702402399595
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.528487
$ws.Range("N5").Value = 85.585461
$ws.Range("O5").Value = 0.1381240089280516
$ws.Range("P5").Value = 0.1381240089280516
$ws.Range("Q5").Value = 39.61383923157266
$ws.Range("R5").Value = 356.5245530841539
$ws.Range("S5").Value = 0.04802604973334666
$ws.Range("T5").Value = 0.04802604973334666

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Amelx"
$ws.Range("C6").Value = "Lamp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.388571333333333
$ws.Range("H6").Value = 4.165713999999999
$ws.Range("I6").Value = 0.347702402399595
$ws.Range("J6").Value = 0.347702402399595
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 139.2310486666667
$ws.Range("N6").Value = 417.693146
$ws.Range("O6").Value = 0.6741034184216166
$ws.Range("P6").Value = 0.6741034184216166
$ws.Range("Q6").Value = 193.3322428884715
$ws.Range("R6").Value = 1739.990185996244
$ws.Range("S6").Value = 0.2343873780509755
$ws.Range("T6").Value = 0.2343873780509755

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Amelx"
$ws.Range("C7").Value = "Lamp2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.388571333333333
$ws.Range("H7").Value = 4.165713999999999
$ws.Range("I7").Value = 0.347702402399595
$ws.Range("J7").Value = 0.347702402399595
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 38.78302866666667
$ws.Range("N7").Value = 116.349086
$ws.Range("O7").Value = 0.1877725726503318
$ws.Range("P7").Value = 0.1877725726503318
$ws.Range("Q7").Value = 53.85300182637822
$ws.Range("R7").Value = 484.6770164374039
$ws.Range("S7").Value = 0.06528897461527287
$ws.Range("T7").Value = 0.06528897461527287

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Amelx"
$ws.Range("C8").Value = "Lamp2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.552285333333333
$ws.Range("H8").Value = 7.656856
$ws.Range("I8").Value = 0.6390998580382028
$ws.Range("J8").Value = 0.6390998580382028
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 28.528487
$ws.Range("N8").Value = 85.585461
$ws.Range("O8").Value = 0.1381240089280516
$ws.Range("P8").Value = 0.1381240089280516
$ws.Range("Q8").Value = 72.81283895229066
$ws.Range("R8").Value = 655.315550570616
$ws.Range("S8").Value = 0.08827503449758524
$ws.Range("T8").Value = 0.08827503449758524

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Amelx"
$ws.Range("C9").Value = "Lamp2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.552285333333333
$ws.Range("H9").Value = 7.656856
$ws.Range("I9").Value = 0.6390998580382028
$ws.Range("J9").Value = 0.6390998580382028
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 139.2310486666667
$ws.Range("N9").Value = 417.693146
$ws.Range("O9").Value = 0.6741034184216166
$ws.Range("P9").Value = 0.6741034184216166
$ws.Range("Q9").Value = 355.3573634565528
$ws.Range("R9").Value = 3198.216271108976
$ws.Range("S9").Value = 0.4308193990163224
$ws.Range("T9").Value = 0.4308193990163224

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Amelx"
$ws.Range("C10").Value = "Lamp2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.552285333333333
$ws.Range("H10").Value = 7.656856
$ws.Range("I10").Value = 0.6390998580382028
$ws.Range("J10").Value = 0.6390998580382028
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 38.78302866666667
$ws.Range("N10").Value = 116.349086
$ws.Range("O10").Value = 0.1877725726503318
$ws.Range("P10").Value = 0.1877725726503318
$ws.Range("Q10").Value = 98.98535524817956
$ws.Range("R10").Value = 890.868197233616
$ws.Range("S10").Value = 0.1200054245242952
$ws.Range("T10").Value = 0.1200054245242952
